# Titular - cargado.xlsx : carga egreso, infractorcontraventor y titular
#
# The "idtipodni" column (L) used to hold a raw numeric code (1 or 2).
# It is converted to hold the actual text label instead:
#   1 -> "DNI"
#   2 -> "LE"
# The sheet view is also scrolled/selected to cell O83.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 101) { $lastRow = 101 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 12)
    $v = $cell.Value2
    if ($v -eq 1) {
        $cell.Value = "DNI"
    } elseif ($v -eq 2) {
        $cell.Value = "LE"
    }
}

# Update the window/sheet view: scroll so column B is the left-most visible
# column, and move the active selection to O83.
$ws.Range("O83").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
